$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 207.9
$ws.Cells.Item(9, 10).Value = 224.83333
$ws.Cells.Item(9, 12).Value = 224.83333
$ws.Cells.Item(9, 14).Value = -562.8333299999999
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(18, 8).Value = 684
$ws.Cells.Item(18, 9).Value = 698.8
$ws.Cells.Item(18, 10).Value = 647
$ws.Cells.Item(18, 11).Value = 698.8
$ws.Cells.Item(18, 12).Value = 647
$ws.Cells.Item(18, 13).Value = -414.8
$ws.Cells.Item(18, 14).Value = -1215
$ws.Cells.Item(19, 8).Value = 1184.6364
$ws.Cells.Item(19, 9).Value = 855.8
$ws.Cells.Item(19, 10).Value = 1458.6666
$ws.Cells.Item(19, 11).Value = 855.8
$ws.Cells.Item(19, 12).Value = 1458.6666
$ws.Cells.Item(19, 13).Value = -680.8
$ws.Cells.Item(19, 14).Value = -1808.6666
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).ClearContents()
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 14).ClearContents()
$ws.Cells.Item(28, 8).Value = 3842.6667
$ws.Cells.Item(28, 9).Value = 3842.6667
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 3842.6667
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = -3357.6667
$ws.Cells.Item(28, 14).ClearContents()
$ws.Cells.Item(41, 8).Value = 829.3333
$ws.Cells.Item(41, 9).Value = 829.3333
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 829.3333
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).Value = -389.3333
$ws.Cells.Item(41, 14).ClearContents()
$ws.Cells.Item(53, 8).Value = 95672.71000000001
$ws.Cells.Item(53, 9).Value = 377.6
$ws.Cells.Item(53, 11).Value = 377.6
$ws.Cells.Item(53, 13).Value = 259.4
$ws.Cells.Item(62, 8).Value = 3745.8
$ws.Cells.Item(62, 9).Value = 3975
$ws.Cells.Item(62, 11).Value = 3975
$ws.Cells.Item(62, 13).Value = -3351
$ws.Cells.Item(65, 8).Value = 3745.8
$ws.Cells.Item(65, 9).Value = 3975
$ws.Cells.Item(65, 11).Value = 19875
$ws.Cells.Item(65, 13).Value = -16755
$ws.Cells.Item(69, 8).Value = 6500
$ws.Cells.Item(69, 9).Value = 6500
$ws.Cells.Item(69, 11).Value = 19500
$ws.Cells.Item(69, 13).Value = -18626
$ws.Cells.Item(72, 8).Value = 6500
$ws.Cells.Item(72, 9).Value = 6500
$ws.Cells.Item(72, 11).Value = 58500
$ws.Cells.Item(72, 13).Value = -54132
$ws.Cells.Item(76, 8).Value = 16785.428
$ws.Cells.Item(76, 9).Value = 12138.4
$ws.Cells.Item(76, 10).Value = 17413.406
$ws.Cells.Item(76, 11).Value = 12138.4
$ws.Cells.Item(76, 12).Value = 17413.406
$ws.Cells.Item(76, 13).Value = -11823.4
$ws.Cells.Item(76, 14).Value = -18043.406
$ws.Cells.Item(79, 8).Value = 16785.428
$ws.Cells.Item(79, 9).Value = 12138.4
$ws.Cells.Item(79, 10).Value = 17413.406
$ws.Cells.Item(79, 11).Value = 12138.4
$ws.Cells.Item(79, 12).Value = 17413.406
$ws.Cells.Item(79, 13).Value = -11046.4
$ws.Cells.Item(79, 14).Value = -19597.406
$ws.Cells.Item(86, 8).Value = 3545.4
$ws.Cells.Item(86, 9).Value = 3309.3333
$ws.Cells.Item(86, 11).Value = 3309.3333
$ws.Cells.Item(86, 13).Value = -2186.3333
$ws.Cells.Item(89, 8).Value = 3545.4
$ws.Cells.Item(89, 9).Value = 3309.3333
$ws.Cells.Item(89, 11).Value = 16546.6665
$ws.Cells.Item(89, 13).Value = -10930.6665
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).ClearContents()
$ws.Cells.Item(92, 14).ClearContents()
$ws.Cells.Item(96, 8).Value = 1195
$ws.Cells.Item(96, 9).Value = 1331.25
$ws.Cells.Item(96, 10).Value = 650
$ws.Cells.Item(96, 11).Value = 3993.75
$ws.Cells.Item(96, 12).Value = 1950
$ws.Cells.Item(96, 13).Value = -2620.75
$ws.Cells.Item(96, 14).Value = -4696
$ws.Cells.Item(98, 8).Value = 739.8
$ws.Cells.Item(98, 9).Value = 549.75
$ws.Cells.Item(98, 11).Value = 549.75
$ws.Cells.Item(98, 13).Value = 948.25
$ws.Cells.Item(106, 8).Value = 3098
$ws.Cells.Item(106, 9).Value = 3020
$ws.Cells.Item(106, 10).Value = 3800
$ws.Cells.Item(106, 11).Value = 3020
$ws.Cells.Item(106, 12).Value = 3800
$ws.Cells.Item(106, 13).Value = -2389
$ws.Cells.Item(106, 14).Value = -5062
$ws.Cells.Item(107, 8).Value = 3418
$ws.Cells.Item(107, 9).Value = 3418
$ws.Cells.Item(107, 11).Value = 3418
$ws.Cells.Item(107, 13).Value = -1498
$ws.Cells.Item(111, 8).Value = 3435.0667
$ws.Cells.Item(111, 9).Value = 1658
$ws.Cells.Item(111, 10).Value = 4323.6
$ws.Cells.Item(111, 11).Value = 4974
$ws.Cells.Item(111, 12).Value = 12970.8
$ws.Cells.Item(111, 13).Value = -1907
$ws.Cells.Item(111, 14).Value = -19104.8
$ws.Cells.Item(112, 8).Value = 3761.2666
$ws.Cells.Item(112, 10).Value = 4672
$ws.Cells.Item(112, 12).Value = 14016
$ws.Cells.Item(112, 14).Value = -16232
$ws.Cells.Item(122, 8).Value = 739.8
$ws.Cells.Item(122, 9).Value = 549.75
$ws.Cells.Item(122, 11).Value = 1649.25
$ws.Cells.Item(122, 13).Value = 800.75
$ws.Cells.Item(132, 8).Value = 4922.2593
$ws.Cells.Item(132, 9).Value = 6119.25
$ws.Cells.Item(132, 11).Value = 18357.75
$ws.Cells.Item(132, 13).Value = -15827.75
$ws.Cells.Item(135, 8).Value = 850.86365
$ws.Cells.Item(135, 9).Value = 733.5
$ws.Cells.Item(135, 11).Value = 6601.5
$ws.Cells.Item(135, 13).Value = -4066.5
$ws.Cells.Item(138, 8).Value = 2098.9473
$ws.Cells.Item(138, 9).Value = 1317.4482
$ws.Cells.Item(138, 10).Value = 4617.1113
$ws.Cells.Item(138, 11).Value = 3952.3446
$ws.Cells.Item(138, 12).Value = 13851.3339
$ws.Cells.Item(138, 13).Value = 1187.6554
$ws.Cells.Item(138, 14).Value = -24131.3339
$ws.Cells.Item(141, 8).Value = 1690.5161
$ws.Cells.Item(141, 9).Value = 1690.5161
$ws.Cells.Item(141, 11).Value = 5071.5483
$ws.Cells.Item(141, 13).Value = 108.4516999999996

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 4896.385
$ws.Cells.Item(2, 9).Value = 4126.25
$ws.Cells.Item(2, 11).Value = 4126.25
$ws.Cells.Item(2, 13).Value = -4013.25
$ws.Cells.Item(32, 8).Value = 3835.7527
$ws.Cells.Item(32, 9).Value = 3176.6853
$ws.Cells.Item(32, 11).Value = 3176.6853
$ws.Cells.Item(32, 13).Value = -2889.6853
$ws.Cells.Item(45, 8).Value = 57042.46
$ws.Cells.Item(45, 9).Value = 73616.82000000001
$ws.Cells.Item(45, 10).Value = 5477.778
$ws.Cells.Item(45, 11).Value = 73616.82000000001
$ws.Cells.Item(45, 12).Value = 5477.778
$ws.Cells.Item(45, 13).Value = -73239.82000000001
$ws.Cells.Item(45, 14).Value = -6231.778
$ws.Cells.Item(97, 8).Value = 1278.5264
$ws.Cells.Item(97, 9).Value = 1246
$ws.Cells.Item(97, 11).Value = 1246
$ws.Cells.Item(97, 13).Value = -750
$ws.Cells.Item(102, 8).Value = 4262
$ws.Cells.Item(102, 9).Value = 2934.2222
$ws.Cells.Item(102, 10).Value = 7249.5
$ws.Cells.Item(102, 11).Value = 2934.2222
$ws.Cells.Item(102, 12).Value = 7249.5
$ws.Cells.Item(102, 13).Value = -1312.2222
$ws.Cells.Item(102, 14).Value = -10493.5
$ws.Cells.Item(110, 8).Value = 5334.7085
$ws.Cells.Item(110, 9).Value = 4301.65
$ws.Cells.Item(110, 11).Value = 4301.65
$ws.Cells.Item(110, 13).Value = -2256.65
$ws.Cells.Item(116, 8).Value = 4896.385
$ws.Cells.Item(116, 9).Value = 4126.25
$ws.Cells.Item(116, 11).Value = 4126.25
$ws.Cells.Item(116, 13).Value = -1832.25
$ws.Cells.Item(122, 8).Value = 1282.25
$ws.Cells.Item(122, 9).Value = 1282.25
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 3846.75
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -1396.75
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 3628.4285
$ws.Cells.Item(132, 9).Value = 3750.963
$ws.Cells.Item(132, 11).Value = 11252.889
$ws.Cells.Item(132, 13).Value = -8722.889000000001

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(2, 8).Value = 69600
$ws.Cells.Item(2, 10).Value = 69600
$ws.Cells.Item(2, 12).Value = 69600
$ws.Cells.Item(2, 14).Value = -69826
$ws.Cells.Item(3, 8).Value = 4896.385
$ws.Cells.Item(3, 9).Value = 4126.25
$ws.Cells.Item(3, 11).Value = 4126.25
$ws.Cells.Item(3, 13).Value = -4012.25
$ws.Cells.Item(47, 8).Value = 200000
$ws.Cells.Item(47, 10).Value = 200000
$ws.Cells.Item(47, 12).Value = 200000
$ws.Cells.Item(47, 14).Value = -201040
$ws.Cells.Item(80, 8).Value = 500
$ws.Cells.Item(80, 9).Value = 268.5
$ws.Cells.Item(80, 10).Value = 654.3333
$ws.Cells.Item(80, 11).Value = 268.5
$ws.Cells.Item(80, 12).Value = 654.3333
$ws.Cells.Item(80, 13).Value = 729.5
$ws.Cells.Item(80, 14).Value = -2650.3333
$ws.Cells.Item(83, 8).Value = 500
$ws.Cells.Item(83, 9).Value = 268.5
$ws.Cells.Item(83, 10).Value = 654.3333
$ws.Cells.Item(83, 11).Value = 1342.5
$ws.Cells.Item(83, 12).Value = 3271.6665
$ws.Cells.Item(83, 13).Value = 3649.5
$ws.Cells.Item(83, 14).Value = -13255.6665
$ws.Cells.Item(94, 8).Value = 2316.2563
$ws.Cells.Item(94, 9).Value = 2100.25
$ws.Cells.Item(94, 11).Value = 2100.25
$ws.Cells.Item(94, 13).Value = -1649.25
$ws.Cells.Item(105, 8).Value = 1850
$ws.Cells.Item(105, 9).Value = 1850
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 1850
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = -103
$ws.Cells.Item(105, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 98990
$ws.Cells.Item(126, 10).Value = 98990
$ws.Cells.Item(126, 12).Value = 98990
$ws.Cells.Item(126, 14).Value = -108870
$ws.Cells.Item(132, 8).Value = 100000
$ws.Cells.Item(132, 10).Value = 100000
$ws.Cells.Item(132, 12).Value = 100000
$ws.Cells.Item(132, 14).Value = -110120
$ws.Cells.Item(134, 8).Value = 3204.2292
$ws.Cells.Item(134, 9).Value = 3204.2292
$ws.Cells.Item(134, 11).Value = 9612.687600000001
$ws.Cells.Item(134, 13).Value = -7077.687600000001

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3370.1538
$ws.Cells.Item(16, 9).Value = 2632.5881
$ws.Cells.Item(16, 11).Value = 2632.5881
$ws.Cells.Item(16, 13).Value = -2345.5881
$ws.Cells.Item(22, 8).Value = 745.44446
$ws.Cells.Item(22, 9).Value = 385.33334
$ws.Cells.Item(22, 11).Value = 385.33334
$ws.Cells.Item(22, 13).Value = -35.33334000000002
$ws.Cells.Item(31, 8).Value = 6068.857
$ws.Cells.Item(31, 9).Value = 3721.2856
$ws.Cells.Item(31, 10).Value = 7242.643
$ws.Cells.Item(31, 11).Value = 3721.2856
$ws.Cells.Item(31, 12).Value = 7242.643
$ws.Cells.Item(31, 13).Value = -3426.2856
$ws.Cells.Item(31, 14).Value = -7832.643
$ws.Cells.Item(34, 8).Value = 6068.857
$ws.Cells.Item(34, 9).Value = 3721.2856
$ws.Cells.Item(34, 10).Value = 7242.643
$ws.Cells.Item(34, 11).Value = 3721.2856
$ws.Cells.Item(34, 12).Value = 7242.643
$ws.Cells.Item(34, 13).Value = -3519.2856
$ws.Cells.Item(34, 14).Value = -7646.643
$ws.Cells.Item(58, 8).Value = 6472.773
$ws.Cells.Item(58, 9).Value = 3900.1
$ws.Cells.Item(58, 11).Value = 3900.1
$ws.Cells.Item(58, 13).Value = -3697.1
$ws.Cells.Item(62, 8).Value = 8917.6
$ws.Cells.Item(62, 9).Value = 8994.5
$ws.Cells.Item(62, 10).Value = 8866.333000000001
$ws.Cells.Item(62, 11).Value = 8994.5
$ws.Cells.Item(62, 12).Value = 8866.333000000001
$ws.Cells.Item(62, 13).Value = -8370.5
$ws.Cells.Item(62, 14).Value = -10114.333
$ws.Cells.Item(65, 8).Value = 8917.6
$ws.Cells.Item(65, 9).Value = 8994.5
$ws.Cells.Item(65, 10).Value = 8866.333000000001
$ws.Cells.Item(65, 11).Value = 44972.5
$ws.Cells.Item(65, 12).Value = 44331.665
$ws.Cells.Item(65, 13).Value = -41852.5
$ws.Cells.Item(65, 14).Value = -50571.665
$ws.Cells.Item(94, 8).Value = 2649.3333
$ws.Cells.Item(94, 10).Value = 2974
$ws.Cells.Item(94, 12).Value = 2974
$ws.Cells.Item(94, 14).Value = -3876
$ws.Cells.Item(97, 8).Value = 39849
$ws.Cells.Item(97, 10).Value = 39849
$ws.Cells.Item(97, 12).Value = 39849
$ws.Cells.Item(97, 14).Value = -41831
$ws.Cells.Item(99, 8).Value = 5917.4736
$ws.Cells.Item(99, 10).Value = 5493.6665
$ws.Cells.Item(99, 12).Value = 5493.6665
$ws.Cells.Item(99, 14).Value = -8489.666499999999
$ws.Cells.Item(113, 8).Value = 3370.1538
$ws.Cells.Item(113, 9).Value = 2632.5881
$ws.Cells.Item(113, 11).Value = 2632.5881
$ws.Cells.Item(113, 13).Value = -462.5880999999999
$ws.Cells.Item(126, 8).Value = 5917.4736
$ws.Cells.Item(126, 10).Value = 5493.6665
$ws.Cells.Item(126, 12).Value = 16480.9995
$ws.Cells.Item(126, 14).Value = -21420.9995
$ws.Cells.Item(132, 8).Value = 3723.875
$ws.Cells.Item(132, 9).Value = 4113.9565
$ws.Cells.Item(132, 10).Value = 2727
$ws.Cells.Item(132, 11).Value = 12341.8695
$ws.Cells.Item(132, 12).Value = 8181
$ws.Cells.Item(132, 13).Value = -9811.869500000001
$ws.Cells.Item(132, 14).Value = -13241
$ws.Cells.Item(134, 8).Value = 4983.7036
$ws.Cells.Item(134, 9).Value = 3952.7
$ws.Cells.Item(134, 11).Value = 11858.1
$ws.Cells.Item(134, 13).Value = -9323.099999999999
$ws.Cells.Item(136, 8).Value = 6472.773
$ws.Cells.Item(136, 9).Value = 3900.1
$ws.Cells.Item(136, 11).Value = 11700.3
$ws.Cells.Item(136, 13).Value = -9150.299999999999
$ws.Cells.Item(141, 8).Value = 37849.5
$ws.Cells.Item(141, 10).Value = 37849.5
$ws.Cells.Item(141, 12).Value = 37849.5
$ws.Cells.Item(141, 14).Value = -48209.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2618.2632
$ws.Cells.Item(5, 9).Value = 5701.3335
$ws.Cells.Item(5, 10).Value = 2040.1875
$ws.Cells.Item(5, 11).Value = 17104.0005
$ws.Cells.Item(5, 12).Value = 6120.5625
$ws.Cells.Item(5, 13).Value = -16992.0005
$ws.Cells.Item(5, 14).Value = -6344.5625
$ws.Cells.Item(26, 8).Value = 252800.25
$ws.Cells.Item(26, 10).Value = 3732.8333
$ws.Cells.Item(26, 12).Value = 11198.4999
$ws.Cells.Item(26, 14).Value = -11774.4999
$ws.Cells.Item(51, 8).Value = 60073.8
$ws.Cells.Item(51, 9).Value = 28937
$ws.Cells.Item(51, 11).Value = 86811
$ws.Cells.Item(51, 13).Value = -86351
$ws.Cells.Item(56, 8).Value = 53500
$ws.Cells.Item(56, 9).Value = 53500
$ws.Cells.Item(56, 11).Value = 53500
$ws.Cells.Item(56, 13).Value = -52970
$ws.Cells.Item(131, 8).Value = 21741448
$ws.Cells.Item(131, 9).Value = 250001000
$ws.Cells.Item(131, 10).Value = 2443.6667
$ws.Cells.Item(131, 11).Value = 750003000
$ws.Cells.Item(131, 12).Value = 7331.000100000001
$ws.Cells.Item(131, 13).Value = -749997960
$ws.Cells.Item(131, 14).Value = -17411.0001
$ws.Cells.Item(135, 8).Value = 2618.2632
$ws.Cells.Item(135, 9).Value = 5701.3335
$ws.Cells.Item(135, 10).Value = 2040.1875
$ws.Cells.Item(135, 11).Value = 51312.0015
$ws.Cells.Item(135, 12).Value = 18361.6875
$ws.Cells.Item(135, 13).Value = -48777.0015
$ws.Cells.Item(135, 14).Value = -23431.6875
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).ClearContents()
$ws.Cells.Item(137, 14).ClearContents()

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10860.186
$ws.Cells.Item(70, 9).Value = 7695.154
$ws.Cells.Item(70, 10).Value = 13799.143
$ws.Cells.Item(70, 11).Value = 7695.154
$ws.Cells.Item(70, 12).Value = 13799.143
$ws.Cells.Item(70, 13).Value = -7425.154
$ws.Cells.Item(70, 14).Value = -14339.143
$ws.Cells.Item(73, 8).Value = 10860.186
$ws.Cells.Item(73, 9).Value = 7695.154
$ws.Cells.Item(73, 10).Value = 13799.143
$ws.Cells.Item(73, 11).Value = 7695.154
$ws.Cells.Item(73, 12).Value = 13799.143
$ws.Cells.Item(73, 13).Value = -6759.154
$ws.Cells.Item(73, 14).Value = -15671.143
$ws.Cells.Item(97, 8).Value = 2162.5134
$ws.Cells.Item(97, 9).Value = 820.625
$ws.Cells.Item(97, 10).Value = 4639.846
$ws.Cells.Item(97, 11).Value = 820.625
$ws.Cells.Item(97, 12).Value = 4639.846
$ws.Cells.Item(97, 13).Value = -324.625
$ws.Cells.Item(97, 14).Value = -5631.846
$ws.Cells.Item(102, 8).Value = 3607.55
$ws.Cells.Item(102, 9).Value = 2675.0557
$ws.Cells.Item(102, 11).Value = 2675.0557
$ws.Cells.Item(102, 13).Value = -1053.0557
$ws.Cells.Item(122, 8).Value = 10994.833
$ws.Cells.Item(122, 9).Value = 3742.25
$ws.Cells.Item(122, 11).Value = 11226.75
$ws.Cells.Item(122, 13).Value = -8776.75
$ws.Cells.Item(132, 8).Value = 2456.606
$ws.Cells.Item(132, 9).Value = 2394
$ws.Cells.Item(132, 11).Value = 7182
$ws.Cells.Item(132, 13).Value = -4652

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2510.75
$ws.Cells.Item(7, 9).Value = 2569.4285
$ws.Cells.Item(7, 11).Value = 2569.4285
$ws.Cells.Item(7, 13).Value = -2457.4285
$ws.Cells.Item(16, 8).Value = 1175.8182
$ws.Cells.Item(16, 9).Value = 993.7143
$ws.Cells.Item(16, 11).Value = 993.7143
$ws.Cells.Item(16, 13).Value = -823.7143
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 14).ClearContents()
$ws.Cells.Item(40, 8).Value = 4165.7144
$ws.Cells.Item(40, 9).Value = 4165.7144
$ws.Cells.Item(40, 11).Value = 4165.7144
$ws.Cells.Item(40, 13).Value = -4029.7144
$ws.Cells.Item(46, 8).Value = 6799.7856
$ws.Cells.Item(46, 10).Value = 2399
$ws.Cells.Item(46, 12).Value = 2399
$ws.Cells.Item(46, 14).Value = -2775
$ws.Cells.Item(55, 8).Value = 372.25
$ws.Cells.Item(55, 9).Value = 554
$ws.Cells.Item(55, 11).Value = 554
$ws.Cells.Item(55, 13).Value = -381
$ws.Cells.Item(82, 8).Value = 835.7619
$ws.Cells.Item(82, 9).Value = 730.2727
$ws.Cells.Item(82, 10).Value = 951.8
$ws.Cells.Item(82, 11).Value = 730.2727
$ws.Cells.Item(82, 12).Value = 951.8
$ws.Cells.Item(82, 13).Value = -369.2727
$ws.Cells.Item(82, 14).Value = -1673.8
$ws.Cells.Item(85, 8).Value = 835.7619
$ws.Cells.Item(85, 9).Value = 730.2727
$ws.Cells.Item(85, 10).Value = 951.8
$ws.Cells.Item(85, 11).Value = 730.2727
$ws.Cells.Item(85, 12).Value = 951.8
$ws.Cells.Item(85, 13).Value = 517.7273
$ws.Cells.Item(85, 14).Value = -3447.8
$ws.Cells.Item(93, 8).Value = 6839.35
$ws.Cells.Item(93, 9).Value = 1574.9656
$ws.Cells.Item(93, 11).Value = 1574.9656
$ws.Cells.Item(93, 13).Value = -326.9656
$ws.Cells.Item(100, 8).Value = 5037.375
$ws.Cells.Item(100, 9).Value = 4383.1665
$ws.Cells.Item(100, 11).Value = 4383.1665
$ws.Cells.Item(100, 13).Value = -3842.1665
$ws.Cells.Item(122, 8).Value = 3970
$ws.Cells.Item(122, 9).Value = 3266.6667
$ws.Cells.Item(122, 11).Value = 9800.000100000001
$ws.Cells.Item(122, 13).Value = -7350.000100000001
$ws.Cells.Item(126, 8).Value = 2510.75
$ws.Cells.Item(126, 9).Value = 2569.4285
$ws.Cells.Item(126, 11).Value = 7708.2855
$ws.Cells.Item(126, 13).Value = -5238.2855
$ws.Cells.Item(132, 8).Value = 12531.213
$ws.Cells.Item(132, 9).Value = 11429.625
$ws.Cells.Item(132, 11).Value = 34288.875
$ws.Cells.Item(132, 13).Value = -31758.875
$ws.Cells.Item(136, 8).Value = 2451.3845
$ws.Cells.Item(136, 9).Value = 2536.5652
$ws.Cells.Item(136, 11).Value = 7609.6956
$ws.Cells.Item(136, 13).Value = -5059.6956

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 666666
$ws.Cells.Item(14, 10).Value = 666666
$ws.Cells.Item(14, 12).Value = 666666
$ws.Cells.Item(14, 14).Value = -667002
$ws.Cells.Item(49, 8).Value = 15000
$ws.Cells.Item(49, 10).Value = 15000
$ws.Cells.Item(49, 12).Value = 15000
$ws.Cells.Item(49, 14).Value = -15460
$ws.Cells.Item(51, 8).Value = 26666
$ws.Cells.Item(51, 10).Value = 32999.5
$ws.Cells.Item(51, 12).Value = 32999.5
$ws.Cells.Item(51, 14).Value = -34019.5
$ws.Cells.Item(70, 8).Value = 29499.5
$ws.Cells.Item(70, 9).Value = 29499.5
$ws.Cells.Item(70, 11).Value = 29499.5
$ws.Cells.Item(70, 13).Value = -29184.5
$ws.Cells.Item(73, 8).Value = 29499.5
$ws.Cells.Item(73, 9).Value = 29499.5
$ws.Cells.Item(73, 11).Value = 29499.5
$ws.Cells.Item(73, 13).Value = -28407.5
$ws.Cells.Item(100, 8).Value = 463.57895
$ws.Cells.Item(100, 9).Value = 356.75
$ws.Cells.Item(100, 11).Value = 713.5
$ws.Cells.Item(100, 13).Value = -172.5
$ws.Cells.Item(107, 8).Value = 752.8421
$ws.Cells.Item(107, 9).Value = 384.75
$ws.Cells.Item(107, 11).Value = 1154.25
$ws.Cells.Item(107, 13).Value = 765.75
$ws.Cells.Item(122, 8).Value = 3759.88
$ws.Cells.Item(122, 9).Value = 2694.3333
$ws.Cells.Item(122, 11).Value = 8082.999899999999
$ws.Cells.Item(122, 13).Value = -5632.999899999999
$ws.Cells.Item(126, 8).Value = 4490
$ws.Cells.Item(126, 9).Value = 5259.4
$ws.Cells.Item(126, 11).Value = 15778.2
$ws.Cells.Item(126, 13).Value = -13308.2
$ws.Cells.Item(132, 8).Value = 1998.4584
$ws.Cells.Item(132, 9).Value = 2112.6365
$ws.Cells.Item(132, 10).Value = 742.5
$ws.Cells.Item(132, 11).Value = 6337.9095
$ws.Cells.Item(132, 12).Value = 2227.5
$ws.Cells.Item(132, 13).Value = -3807.9095
$ws.Cells.Item(132, 14).Value = -7287.5
$ws.Cells.Item(133, 8).Value = 125342.6
$ws.Cells.Item(133, 9).Value = 46000
$ws.Cells.Item(133, 10).Value = 145178.25
$ws.Cells.Item(133, 11).Value = 46000
$ws.Cells.Item(133, 12).Value = 145178.25
$ws.Cells.Item(133, 13).Value = -40940
$ws.Cells.Item(133, 14).Value = -155298.25
$ws.Cells.Item(136, 8).Value = 3730.7273
$ws.Cells.Item(136, 9).Value = 1896.6774
$ws.Cells.Item(136, 11).Value = 5690.0322
$ws.Cells.Item(136, 13).Value = -3140.0322
